$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.279.43'
$ws.Range("E2").Value = '  +4.03%  '
$ws.Range("D3").Value = '1.703.70'
$ws.Range("E3").Value = '  +1.09%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '''217.72'
$ws.Range("E5").Value = '  +0.79%  '
$ws.Range("E6").Value = '  +0.37%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = '''24.13'
$ws.Range("E8").Value = '  +4.16%  '
$ws.Range("E9").Value = '  +2.09%  '
$ws.Range("D10").Value = '''0.0630'
$ws.Range("E10").Value = '  +0.61%  '
$ws.Range("D11").Value = '''0.0888'
$ws.Range("E11").Value = '  -0.14%  '
$ws.Range("D12").Value = '1.952.07'
$ws.Range("E12").Value = '  +1.44%  '
$ws.Range("D13").Value = '1.708.90'
$ws.Range("E13").Value = '  +1.32%  '
$ws.Range("D14").Value = '''4.19'
$ws.Range("E14").Value = '  -0.06%  '
$ws.Range("D15").Value = '''0.558'
$ws.Range("E15").Value = '  +0.66%  '
$ws.Range("D16").Value = '''67.22'
$ws.Range("E16").Value = '  +0.43%  '
$ws.Range("D17").Value = '''252.13'
$ws.Range("E17").Value = '  +6.76%  '
$ws.Range("D18").Value = '28.275.04'
$ws.Range("E18").Value = '  +4.00%  '
$ws.Range("D19").Value = '0.0₃0745'
$ws.Range("E19").Value = '  +0.12%  '
$ws.Range("D20").Value = '''7.69'
$ws.Range("E20").Value = '  -4.18%  '
$ws.Range("D21").Value = '''1.00'
$ws.Range("E21").Value = '  -0.14%  '
$ws.Range("D22").Value = '''4.55'
$ws.Range("E22").Value = '  -0.42%  '
$ws.Range("D23").Value = '''9.57'
$ws.Range("E23").Value = '  -0.42%  '
$ws.Range("D24").Value = '''2.04'
$ws.Range("E24").Value = '  -1.79%  '
$ws.Range("D25").Value = '''147.53'
$ws.Range("E25").Value = '  +0.02%  '
$ws.Range("D26").Value = '''7.35'
$ws.Range("E26").Value = '  +0.20%  '
$ws.Range("D27").Value = '''16.55'
$ws.Range("E27").Value = '  +0.60%  '
$ws.Range("D28").Value = '''0.113'
$ws.Range("E28").Value = '  +0.24%  '
$ws.Range("E29").Value = '  +0.12%  '
$ws.Range("D30").Value = '''0.0509'
$ws.Range("E30").Value = '  +0.70%  '
$ws.Range("E31").Value = '  +3.20%  '
$ws.Range("E32").Value = '  +0.26%  '
$ws.Range("D33").Value = '1.479.66'
$ws.Range("E33").Value = '  -3.95%  '
$ws.Range("D34").Value = '''3.19'
$ws.Range("E34").Value = '  -1.85%  '
$ws.Range("D35").Value = '''1.62'
$ws.Range("E35").Value = '  -2.78%  '
$ws.Range("D36").Value = '''0.959'
$ws.Range("E36").Value = '  +1.42%  '
$ws.Range("E37").Value = '  +0.19%  '
$ws.Range("E38").Value = '  -1.55%  '
$ws.Range("E39").Value = '  +0.32%  '
$ws.Range("D40").Value = '''1.05'
$ws.Range("E40").Value = '  -1.39%  '
$ws.Range("D41").Value = '''69.24'
$ws.Range("E41").Value = '  +0.11%  '
$ws.Range("E42").Value = '  -0.10%  '
$ws.Range("D43").Value = '''5.62'
$ws.Range("E43").Value = '  -2.32%  '
$ws.Range("D44").Value = '1.855.09'
$ws.Range("E44").Value = '  +1.17%  '
$ws.Range("D45").Value = '''2.25'
$ws.Range("E45").Value = '  -0.23%  '
$ws.Range("E46").Value = '  +1.35%  '
$ws.Range("E47").Value = '  +7.21%  '
$ws.Range("D48").Value = '''89.78'
$ws.Range("E48").Value = '  -0.39%  '
$ws.Range("D49").Value = '0.0₆0107'
$ws.Range("E49").Value = '  -4.38%  '
$ws.Range("E50").Value = '  -0.59%  '
$ws.Range("D51").Value = '''8.03'
$ws.Range("E51").Value = '  -2.66%  '
